$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update cells whose reference already existed with a styled cell; only their text changes
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "8643537 - Fabio Rodolfo Miguel Batista"
$ws.Range("C10").Value = "8643537 - Fabio Rodolfo Miguel Batista"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Duas provas escritas: P1 e P2"
$ws.Range("C19").Value = "Duas provas escritas: P1 e P2"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média das notas obtidas nas duas provas: N1=(P1 + P2)/2"
$ws.Range("C20").Value = "Média das notas obtidas nas duas provas: N1=(P1 + P2)/2"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Uma prova escrita: REC`nMédia das notas N1 e REC:N2=(N1+REC)/2"
$ws.Range("C21").Value = "Uma prova escrita: REC`nMédia das notas N1 e REC:N2=(N1+REC)/2"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B24").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n"
$ws.Range("B25").Value = "LOQ4009 -  Instrumentação na Industria Química  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOQ4009 -  Instrumentação na Industria Química  (Requisito fraco)`n"

# 2) Remove cells that must not exist anymore in the target layout
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# 3) Create brand-new cells: copy formatting from a same-column cell first (so the
#    correct shared cell style is reused instead of a freshly invented one), then set text
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "01/01/2013"
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "01/01/2013"
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "8643537 - Fabio Rodolfo Miguel Batista"
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "8643537 - Fabio Rodolfo Miguel Batista"
$ws.Range("B10").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"
$ws.Range("C10").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"
$excel.CutCopyMode = 0

# 4) Row heights: rows that gain/keep an explicit custom height
$ws.Rows("10:10").RowHeight = 60
$ws.Rows("11:11").RowHeight = 60
$ws.Rows("13:13").RowHeight = 60
$ws.Rows("14:14").RowHeight = 60
$ws.Rows("15:15").RowHeight = 120
$ws.Rows("16:16").RowHeight = 120
$ws.Rows("18:18").RowHeight = 60
$ws.Rows("19:19").RowHeight = 60
$ws.Rows("20:20").RowHeight = 60
$ws.Rows("21:21").RowHeight = 120
$ws.Rows("23:23").RowHeight = 30
$ws.Rows("24:24").RowHeight = 30
$ws.Rows("25:25").RowHeight = 30

# 5) Rows that must lose their custom height entirely (back to sheet default, no ht attribute)
$ws.Rows("17:17").EntireRow.AutoFit()
$ws.Rows("22:22").EntireRow.AutoFit()

# 6) Remove the now-obsolete trailing row (old row 26), which shifts nothing since it is last
$ws.Rows("26:26").Delete()

